$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Backlog sheet: mark stories US01, US02, US04, US05 as selected/"Planned"
# ---------------------------------------------------------------------------
$backlog = $wb.Worksheets.Item("Backlog")

$backlog.Range("A10").Value = 1
$backlog.Range("E10").Value = "Planned"

$backlog.Range("A11").Value = 1
$backlog.Range("E11").Value = "Planned"

$backlog.Range("A18").Value = 1
$backlog.Range("E18").Value = "Planned"

$backlog.Range("A19").Value = 1
$backlog.Range("E19").Value = "Planned"

$backlog.Activate()
$backlog.Range("E12").Select()

# ---------------------------------------------------------------------------
# Sprint1 sheet: add newly-planned stories US01, US02, US04, US05
# ---------------------------------------------------------------------------
$sprint1 = $wb.Worksheets.Item("Sprint1")

$sprint1.Range("A4").Value = "US01"
$sprint1.Range("B4").Value = "Date before current dates"
$sprint1.Range("E4").Value = 100
$sprint1.Range("F4").Value = 120

$sprint1.Range("A5").Value = "US02"
$sprint1.Range("B5").Value = "Birth before marriage"
$sprint1.Range("E5").Value = 100
$sprint1.Range("F5").Value = 120

$sprint1.Range("A8").Value = "US04"
$sprint1.Range("B8").Value = "Marriage before divorce"
$sprint1.Range("E8").Value = 90
$sprint1.Range("F8").Value = 120

$sprint1.Range("A9").Value = "US05"
$sprint1.Range("B9").Value = "Marriage before death"
$sprint1.Range("E9").Value = 90
$sprint1.Range("F9").Value = 120

$sprint1.Activate()
$sprint1.Range("F9").Select()
$excel.ActiveWindow.Zoom = 130

# ---------------------------------------------------------------------------
# Make Team the active/selected sheet (was Stories)
# ---------------------------------------------------------------------------
$team = $wb.Worksheets.Item("Team")
$team.Activate()
